$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1. Row 59 already exists ("Document Repository-default read access" /
#    DOC_REPO / grant read ...). The rule name text gets a small correction
#    and the row becomes the first of a new block of DocumentRepository
#    access-control rules, so update its text in place.
# ---------------------------------------------------------------------------
$ws.Range("B59").Value = "Document Repository -default read access"
$ws.Range("C59").Value = "DOC_REPO"
$ws.Range("G59").Value = "grant read to assignee, co-owner, supervisor, owning group, approver, collaborator, follower, reader, *"

# ---------------------------------------------------------------------------
# 2. Insert 8 new rows (60-67) for the new DocumentRepository rules. Only
#    columns B:G are shifted down so column A (which has no content for
#    these rows) stays untouched/empty.
# ---------------------------------------------------------------------------
$ws.Range("B60:G67").Insert(-4121)
$ws.Range("A60:A67").Clear()

# ---------------------------------------------------------------------------
# 3. Match formatting of the surrounding rule rows by copying the style of
#    representative cells onto the new ranges.
# ---------------------------------------------------------------------------
# Rows 60-63: every column (B:G) uses the plain "rule text" style.
$ws.Range("B59").Copy()
$ws.Range("B60:G63").PasteSpecial(-4122)

# Row 64: column B keeps the plain style, C:G use the "table interior" style.
$ws.Range("B59").Copy()
$ws.Range("B64").PasteSpecial(-4122)
$ws.Range("C59").Copy()
$ws.Range("C64:G64").PasteSpecial(-4122)

# Rows 65-67: B and G use the plain style, C/E/F use the interior style.
$ws.Range("B59").Copy()
$ws.Range("B65:B67").PasteSpecial(-4122)
$ws.Range("G65:G67").PasteSpecial(-4122)
$ws.Range("C59").Copy()
$ws.Range("C65:C67").PasteSpecial(-4122)
$ws.Range("E65:F67").PasteSpecial(-4122)

# D65/D66 use the boolean-formatted style (same as D48), D67 uses the
# boolean-formatted+alignment style (same as D52).
$ws.Range("D48").Copy()
$ws.Range("D65:D66").PasteSpecial(-4122)
$ws.Range("D52").Copy()
$ws.Range("D67").PasteSpecial(-4122)

$ws.Application.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 4. Row heights for the new rows.
# ---------------------------------------------------------------------------
$ws.Rows.Item(60).RowHeight = 45
$ws.Rows.Item(61).RowHeight = 30
$ws.Rows.Item(62).RowHeight = 30
$ws.Rows.Item(63).RowHeight = 30
$ws.Rows.Item(64).RowHeight = 30
$ws.Rows.Item(65).RowHeight = 45
$ws.Rows.Item(66).RowHeight = 45
$ws.Rows.Item(67).RowHeight = 60

# ---------------------------------------------------------------------------
# 5. Values for the new DocumentRepository access-control rules.
# ---------------------------------------------------------------------------
$ws.Range("B60").Value = "DocumentRepository – Anybody can add comments"
$ws.Range("C60").Value = "DOC_REPO"
$ws.Range("G60").Value = "grant addComment to *"

$ws.Range("B61").Value = "DocumentRepository – Lockout No Access Users"
$ws.Range("C61").Value = "DOC_REPO"
$ws.Range("G61").Value = "mandatory deny read to No Access"

$ws.Range("B62").Value = "DocumentRepository – Anybody can add tag"
$ws.Range("C62").Value = "DOC_REPO"
$ws.Range("G62").Value = "grant addTag to *"

$ws.Range("B63").Value = "DocumentRepository – Anybody can subscribe"
$ws.Range("C63").Value = "DOC_REPO"
$ws.Range("G63").Value = "grant subscribe to *"

$ws.Range("B64").Value = "DocumentRepository – Restricted Flag"
$ws.Range("C64").Value = "DOC_REPO"
$ws.Range("D64").Value = "restricted"
$ws.Range("G64").Value = "deny read to *"

$ws.Range("B65").Value = "DocumentRepository – Only participants can add files"
$ws.Range("C65").Value = "DOC_REPO"
$ws.Range("G65").Value = "grant add file to assignee, co-owner, supervisor, owning group, approver, collaborator"

$ws.Range("B66").Value = "DocumentRepository – Only participants can save"
$ws.Range("C66").Value = "DOC_REPO"
$ws.Range("G66").Value = "grant save to assignee, co-owner, supervisor, owning group, approver, collaborator"

$ws.Range("B67").Value = "DocumentRepository – Only participants can upload or replace files"
$ws.Range("C67").Value = "DOC_REPO"
$ws.Range("G67").Value = "grant uploadOrReplaceFile to assignee, co-owner, supervisor, owning group, approver, collaborator, reader"

# ---------------------------------------------------------------------------
# 6. Match the view state recorded for the edited sheet.
# ---------------------------------------------------------------------------
$ws.Range("G67").Select()
$excel.ActiveWindow.ScrollRow = 58

Write-Output "edit complete"
